$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 14 - this shifts existing rows 14:24 down to 15:25,
# carrying all their data/formatting along (matches the diff's row-shift pattern).
$ws.Rows("14:14").Insert()

# Populate the newly inserted row 14 with the new weekly price entry.
$ws.Range("A14").Value = 11
$ws.Range("B14").Value = "Vega Monumental Concepción"
$ws.Range("C14").Value = "Bíobío"
$ws.Range("D14").Value = 44540
$ws.Range("D14").NumberFormat = $ws.Range("D15").NumberFormat
$ws.Range("E14").Value = 8
$ws.Range("F14").Value = 100112026
$ws.Range("G14").Value = "Haba"
$ws.Range("H14").Value = "Sin especificar"
$ws.Range("I14").Value = "Primera"
$ws.Range("J14").Value = 140
$ws.Range("K14").Value = 11000
$ws.Range("L14").Value = 12000
$ws.Range("M14").Value = 11429
$ws.Range("N14").Value = "$/saco 25 kilos"
$ws.Range("O14").Value = "Región del Maule"
$ws.Range("P14").Value = 457
$ws.Range("Q14").Value = 25
$ws.Range("R14").Value = "Hortaliza"
